$wb = $excel.ActiveWorkbook

# --- Add the new ACTOR_AWARD sheet right after DIRECTOR_AWARD ---
$directorAward = $wb.Worksheets.Item("DIRECTOR_AWARD")
$ws = $wb.Worksheets.Add($null, $directorAward)
$ws.Name = "ACTOR_AWARD"

# --- Populate ACTOR_AWARD with the ROLE/award dataset ---
$ws.Cells.Item(1,1).Value = "title"
$ws.Cells.Item(1,2).Value = "production_year"
$ws.Cells.Item(1,3).Value = "description"
$ws.Cells.Item(1,4).Value = "award_name"
$ws.Cells.Item(1,5).Value = "year_of_award"
$ws.Cells.Item(1,6).Value = "category"
$ws.Cells.Item(1,7).Value = "result"
$ws.Cells.Item(2,1).Value = "Traffic"
$ws.Cells.Item(2,2).Value = 2000
$ws.Cells.Item(2,3).Value = "JavierRodriguez"
$ws.Cells.Item(2,4).Value = "Oscar"
$ws.Cells.Item(2,5).Value = 2001
$ws.Cells.Item(2,6).Value = "bestactorinasupportingrole"
$ws.Cells.Item(2,7).Value = "won"
$ws.Cells.Item(3,1).Value = "Gladiator"
$ws.Cells.Item(3,2).Value = 2000
$ws.Cells.Item(3,3).Value = "Maximus"
$ws.Cells.Item(3,4).Value = "Oscar"
$ws.Cells.Item(3,5).Value = 2001
$ws.Cells.Item(3,6).Value = "Bestactorinaleadingrole"
$ws.Cells.Item(3,7).Value = "won"
$ws.Cells.Item(4,1).Value = "LifeisBeautiful"
$ws.Cells.Item(4,2).Value = 1997
$ws.Cells.Item(4,3).Value = "GuidoOrefice"
$ws.Cells.Item(4,4).Value = "Oscar"
$ws.Cells.Item(4,5).Value = 1998
$ws.Cells.Item(4,6).Value = "ActorinaLeadingRole"
$ws.Cells.Item(4,7).Value = "Won"
$ws.Cells.Item(5,1).Value = "Affliction"
$ws.Cells.Item(5,2).Value = 1997
$ws.Cells.Item(5,3).Value = "GlenWhitehouse"
$ws.Cells.Item(5,4).Value = "Oscar"
$ws.Cells.Item(5,5).Value = 1998
$ws.Cells.Item(5,6).Value = "ActorinasupportingRole"
$ws.Cells.Item(5,7).Value = "won"
$ws.Cells.Item(6,1).Value = "ShakespeareinLove"
$ws.Cells.Item(6,2).Value = 1998
$ws.Cells.Item(6,3).Value = "QueenElizabeth"
$ws.Cells.Item(6,4).Value = "Oscar"
$ws.Cells.Item(6,5).Value = 1998
$ws.Cells.Item(6,6).Value = "ActressinasupportingRole"
$ws.Cells.Item(6,7).Value = "won"
$ws.Cells.Item(7,1).Value = "BoysDontCry"
$ws.Cells.Item(7,2).Value = 1999
$ws.Cells.Item(7,3).Value = "BrandonTeena"
$ws.Cells.Item(7,4).Value = "Oscar"
$ws.Cells.Item(7,5).Value = 1999
$ws.Cells.Item(7,6).Value = "ActressinaLeadingRole"
$ws.Cells.Item(7,7).Value = "won"
$ws.Cells.Item(8,1).Value = "AmericanBeauty"
$ws.Cells.Item(8,2).Value = 1999
$ws.Cells.Item(8,3).Value = "LesterBurnham"
$ws.Cells.Item(8,4).Value = "Oscar"
$ws.Cells.Item(8,5).Value = 1999
$ws.Cells.Item(8,6).Value = "ActorinaLeadingRole"
$ws.Cells.Item(8,7).Value = "Won"
$ws.Cells.Item(9,1).Value = "TheCiderHouseRules"
$ws.Cells.Item(9,2).Value = 1999
$ws.Cells.Item(9,3).Value = "DrWilburLarch"
$ws.Cells.Item(9,4).Value = "Oscar"
$ws.Cells.Item(9,5).Value = 1999
$ws.Cells.Item(9,6).Value = "ActorinasupportingRole"
$ws.Cells.Item(9,7).Value = "Won"
$ws.Cells.Item(10,1).Value = "ToplessWomenTalkAboutTheirLives"
$ws.Cells.Item(10,2).Value = 1997
$ws.Cells.Item(10,3).Value = "Neil"
$ws.Cells.Item(10,4).Value = "NewZealandFilmandTVAwards"
$ws.Cells.Item(10,5).Value = 1999
$ws.Cells.Item(10,6).Value = "BestActor"
$ws.Cells.Item(10,7).Value = "won"
$ws.Cells.Item(11,1).Value = "ThePiano"
$ws.Cells.Item(11,2).Value = 1993
$ws.Cells.Item(11,3).Value = "AdaMcGrath"
$ws.Cells.Item(11,4).Value = "Oscar"
$ws.Cells.Item(11,5).Value = 1994
$ws.Cells.Item(11,6).Value = "BestActress"
$ws.Cells.Item(11,7).Value = "won"
$ws.Cells.Item(12,1).Value = "ThePiano"
$ws.Cells.Item(12,2).Value = 1993
$ws.Cells.Item(12,3).Value = "FloraMcGrath"
$ws.Cells.Item(12,4).Value = "Oscar"
$ws.Cells.Item(12,5).Value = 1994
$ws.Cells.Item(12,6).Value = "BestsupportingActress"
$ws.Cells.Item(12,7).Value = "won"
$ws.Cells.Item(13,1).Value = "ThePiano"
$ws.Cells.Item(13,2).Value = 1993
$ws.Cells.Item(13,3).Value = "AdaMcGrath"
$ws.Cells.Item(13,4).Value = "AFIAward"
$ws.Cells.Item(13,5).Value = 1993
$ws.Cells.Item(13,6).Value = "BestActressinaleadrole"
$ws.Cells.Item(13,7).Value = "won"
$ws.Cells.Item(14,1).Value = "ThePiano"
$ws.Cells.Item(14,2).Value = 1993
$ws.Cells.Item(14,3).Value = "GeorgeBaines"
$ws.Cells.Item(14,4).Value = "AFIAward"
$ws.Cells.Item(14,5).Value = 1993
$ws.Cells.Item(14,6).Value = "BestActorinaleadrole"
$ws.Cells.Item(14,7).Value = "won"
$ws.Cells.Item(15,1).Value = "ThePiano"
$ws.Cells.Item(15,2).Value = 1993
$ws.Cells.Item(15,3).Value = "AdaMcGrath"
$ws.Cells.Item(15,4).Value = "BAFTAFilmAward"
$ws.Cells.Item(15,5).Value = 1994
$ws.Cells.Item(15,6).Value = "BestActressinaleadrole"
$ws.Cells.Item(15,7).Value = "won"
$ws.Cells.Item(16,1).Value = "StrictlyBallroom"
$ws.Cells.Item(16,2).Value = 1992
$ws.Cells.Item(16,3).Value = "ShirleyHastings"
$ws.Cells.Item(16,4).Value = "AFIAward"
$ws.Cells.Item(16,5).Value = 1992
$ws.Cells.Item(16,6).Value = "BestActressinasupportingrole"
$ws.Cells.Item(16,7).Value = "won"
$ws.Cells.Item(17,1).Value = "StrictlyBallroom"
$ws.Cells.Item(17,2).Value = 1992
$ws.Cells.Item(17,3).Value = "DougHastings"
$ws.Cells.Item(17,4).Value = "AFIAward"
$ws.Cells.Item(17,5).Value = 1992
$ws.Cells.Item(17,6).Value = "BestActorinasupportingrole"
$ws.Cells.Item(17,7).Value = "Won"
$ws.Cells.Item(18,1).Value = "Traffic"
$ws.Cells.Item(18,2).Value = 2000
$ws.Cells.Item(18,3).Value = "JavierRodriguez"
$ws.Cells.Item(18,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(18,5).Value = 2001
$ws.Cells.Item(18,6).Value = "BestActorinasupportingrole"
$ws.Cells.Item(18,7).Value = "won"
$ws.Cells.Item(19,1).Value = "Traffic"
$ws.Cells.Item(19,2).Value = 2000
$ws.Cells.Item(19,3).Value = "JavierRodriguez"
$ws.Cells.Item(19,4).Value = "SilverBerlinBear"
$ws.Cells.Item(19,5).Value = 2001
$ws.Cells.Item(19,6).Value = "BestActor"
$ws.Cells.Item(19,7).Value = "won"
$ws.Cells.Item(20,1).Value = "Psycho"
$ws.Cells.Item(20,2).Value = 1960
$ws.Cells.Item(20,3).Value = "MarionCrane"
$ws.Cells.Item(20,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(20,5).Value = 1961
$ws.Cells.Item(20,6).Value = "BestSupportingActress"
$ws.Cells.Item(20,7).Value = "Won"
$ws.Cells.Item(21,1).Value = "TwelveMonkeys"
$ws.Cells.Item(21,2).Value = 1995
$ws.Cells.Item(21,3).Value = "JeffreyGoines"
$ws.Cells.Item(21,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(21,5).Value = 1996
$ws.Cells.Item(21,6).Value = "BestSupportingActor"
$ws.Cells.Item(21,7).Value = "won"
$ws.Cells.Item(22,1).Value = "Alice"
$ws.Cells.Item(22,2).Value = 1990
$ws.Cells.Item(22,3).Value = "AliceTate"
$ws.Cells.Item(22,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(22,5).Value = 1991
$ws.Cells.Item(22,6).Value = "BestPerformancebyanActress"
$ws.Cells.Item(22,7).Value = "nominated"
$ws.Cells.Item(23,1).Value = "Chaplin"
$ws.Cells.Item(23,2).Value = 1992
$ws.Cells.Item(23,3).Value = "CharlieChaplin"
$ws.Cells.Item(23,4).Value = "Oscar"
$ws.Cells.Item(23,5).Value = 1993
$ws.Cells.Item(23,6).Value = "BestActorinaLeadingRole"
$ws.Cells.Item(23,7).Value = "nominated"
$ws.Cells.Item(24,1).Value = "Chaplin"
$ws.Cells.Item(24,2).Value = 1992
$ws.Cells.Item(24,3).Value = "CharlieChaplin"
$ws.Cells.Item(24,4).Value = "BAFTAFilmAward"
$ws.Cells.Item(24,5).Value = 1993
$ws.Cells.Item(24,6).Value = "BestActor"
$ws.Cells.Item(24,7).Value = "won"
$ws.Cells.Item(25,1).Value = "Chaplin"
$ws.Cells.Item(25,2).Value = 1992
$ws.Cells.Item(25,3).Value = "CharlieChaplin"
$ws.Cells.Item(25,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(25,5).Value = 1993
$ws.Cells.Item(25,6).Value = "BestPerformancebyanActor"
$ws.Cells.Item(25,7).Value = "nominated"
$ws.Cells.Item(26,1).Value = "Chaplin"
$ws.Cells.Item(26,2).Value = 1992
$ws.Cells.Item(26,3).Value = "HannahChaplin"
$ws.Cells.Item(26,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(26,5).Value = 1993
$ws.Cells.Item(26,6).Value = "BestPerformancebyanActress"
$ws.Cells.Item(26,7).Value = "nominated"
$ws.Cells.Item(27,1).Value = "Chaplin"
$ws.Cells.Item(27,2).Value = 1992
$ws.Cells.Item(27,3).Value = "CharlieChaplin"
$ws.Cells.Item(27,4).Value = "ALFSAward"
$ws.Cells.Item(27,5).Value = 1993
$ws.Cells.Item(27,6).Value = "ActoroftheYear"
$ws.Cells.Item(27,7).Value = "won"
$ws.Cells.Item(28,1).Value = "Fearless"
$ws.Cells.Item(28,2).Value = 1993
$ws.Cells.Item(28,3).Value = "CarlaRodrigo"
$ws.Cells.Item(28,4).Value = "Oscar"
$ws.Cells.Item(28,5).Value = 1994
$ws.Cells.Item(28,6).Value = "BestActressinaSupportingRole"
$ws.Cells.Item(28,7).Value = "nominated"
$ws.Cells.Item(29,1).Value = "Fearless"
$ws.Cells.Item(29,2).Value = 1993
$ws.Cells.Item(29,3).Value = "CarlaRodrigo"
$ws.Cells.Item(29,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(29,5).Value = 1994
$ws.Cells.Item(29,6).Value = "BestActressinaSupportingRole"
$ws.Cells.Item(29,7).Value = "nominated"
$ws.Cells.Item(30,1).Value = "CyranodeBergerac"
$ws.Cells.Item(30,2).Value = 1990
$ws.Cells.Item(30,3).Value = "CyranoDeBergerac"
$ws.Cells.Item(30,4).Value = "Oscar"
$ws.Cells.Item(30,5).Value = 1991
$ws.Cells.Item(30,6).Value = "BestActorinaLeadingRole"
$ws.Cells.Item(30,7).Value = "nominated"
$ws.Cells.Item(31,1).Value = "CyranodeBergerac"
$ws.Cells.Item(31,2).Value = 1990
$ws.Cells.Item(31,3).Value = "CyranoDeBergerac"
$ws.Cells.Item(31,4).Value = "BAFTAFilmAward"
$ws.Cells.Item(31,5).Value = 1992
$ws.Cells.Item(31,6).Value = "BestActor"
$ws.Cells.Item(31,7).Value = "nominated"
$ws.Cells.Item(32,1).Value = "CyranodeBergerac"
$ws.Cells.Item(32,2).Value = 1990
$ws.Cells.Item(32,3).Value = "CyranoDeBergerac"
$ws.Cells.Item(32,4).Value = "ALFSAward"
$ws.Cells.Item(32,5).Value = 1992
$ws.Cells.Item(32,6).Value = "ActorOFtheYear"
$ws.Cells.Item(32,7).Value = "won"
$ws.Cells.Item(33,1).Value = "ManhattanMurderMystery"
$ws.Cells.Item(33,2).Value = 1993
$ws.Cells.Item(33,3).Value = "CarolLipton"
$ws.Cells.Item(33,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(33,5).Value = 1994
$ws.Cells.Item(33,6).Value = "BestPerformancebyanActress"
$ws.Cells.Item(33,7).Value = "nominated"
$ws.Cells.Item(34,1).Value = "BennyandJoon"
$ws.Cells.Item(34,2).Value = 1993
$ws.Cells.Item(34,3).Value = "Sam"
$ws.Cells.Item(34,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(34,5).Value = 1994
$ws.Cells.Item(34,6).Value = "BestPerformancebyanActor"
$ws.Cells.Item(34,7).Value = "nominated"
$ws.Cells.Item(35,1).Value = "SixDegreesofSeparation"
$ws.Cells.Item(35,2).Value = 1993
$ws.Cells.Item(35,3).Value = "OuisaKittredge"
$ws.Cells.Item(35,4).Value = "Oscar"
$ws.Cells.Item(35,5).Value = 1994
$ws.Cells.Item(35,6).Value = "BestActressinaLeadingRole"
$ws.Cells.Item(35,7).Value = "nominated"
$ws.Cells.Item(36,1).Value = "SixDegreesofSeparation"
$ws.Cells.Item(36,2).Value = 1993
$ws.Cells.Item(36,3).Value = "OuisaKittredge"
$ws.Cells.Item(36,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(36,5).Value = 1994
$ws.Cells.Item(36,6).Value = "BestPerformancebyanActress"
$ws.Cells.Item(36,7).Value = "nominated"
$ws.Cells.Item(37,1).Value = "BawangBieJi"
$ws.Cells.Item(37,2).Value = 1993
$ws.Cells.Item(37,3).Value = "Juxian"
$ws.Cells.Item(37,4).Value = "NYFCCAward"
$ws.Cells.Item(37,5).Value = 1993
$ws.Cells.Item(37,6).Value = "BestSupportingActress"
$ws.Cells.Item(37,7).Value = "won"
$ws.Cells.Item(38,1).Value = "IntheLineofFire"
$ws.Cells.Item(38,2).Value = 1993
$ws.Cells.Item(38,3).Value = "MitchLeary"
$ws.Cells.Item(38,4).Value = "Oscar"
$ws.Cells.Item(38,5).Value = 1994
$ws.Cells.Item(38,6).Value = "BestActorinaSupportingRole"
$ws.Cells.Item(38,7).Value = "nominated"
$ws.Cells.Item(39,1).Value = "IntheLineofFire"
$ws.Cells.Item(39,2).Value = 1993
$ws.Cells.Item(39,3).Value = "MitchLeary"
$ws.Cells.Item(39,4).Value = "BAFTAFilmAward"
$ws.Cells.Item(39,5).Value = 1994
$ws.Cells.Item(39,6).Value = "BestActorSupporting"
$ws.Cells.Item(39,7).Value = "nominated"
$ws.Cells.Item(40,1).Value = "IntheLineofFire"
$ws.Cells.Item(40,2).Value = 1993
$ws.Cells.Item(40,3).Value = "MitchLeary"
$ws.Cells.Item(40,4).Value = "GoldenGlobeAwards"
$ws.Cells.Item(40,5).Value = 1994
$ws.Cells.Item(40,6).Value = "BestPerformancebyanActor"
$ws.Cells.Item(40,7).Value = "nominated"

# --- Set the selection/scroll state on the new sheet ---
$ws.Range("B13").Select()
$excel.ActiveWindow.ScrollRow = 6

# --- Update selections on MOVIE and PERSON sheets ---
$movie = $wb.Worksheets.Item("MOVIE")
$movie.Range("E2").Select()

$person = $wb.Worksheets.Item("PERSON")
$person.Range("A2").Select()

# --- Make AWARD the active sheet/tab ---
$award = $wb.Worksheets.Item("AWARD")
$award.Activate()
